$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('150BB4224N', 'Bag Paper - Baguette', '2', '$115.37', '$230.74'),
    @('245CCGR2518', 'Cake Board - Full Sheet (SO)', '2', '$48.31', '$96.62'),
    @('245CCGR1914', 'Cake Board - 1/2 Sheet', '2', '$31.74', '$63.48'),
    @('77031906', 'Loaf Pan - Large Rectangle (paper)', '1', '$134.99', '$134.99'),
    @('760SOUP32PBL', 'Lid - Soup (32oz)', '1', '$37.99', '$37.99'),
    @('760SOUP32MB', 'Container - Soup (32oz)', '1', '$59.99', '$59.99'),
    @('500L4B', 'Lid Espresso - 4oz', '1', '$33.99', '$33.99'),
    @('5004CAFE', 'Cup - Espresso (4oz)', '1', '$32.99', '$32.99'),
    @('43312MINCUP125', 'Container - Mini Muffin (12 Pack)', '1', '$38.49', '$38.49')
)

$startRow = 9
$numRows = $data.Count
$endRow = $startRow + $numRows - 1

# Force the whole new block to be stored as text, matching the rest of the sheet
$ws.Range("A$startRow`:E$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $numRows; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
    $ws.Cells.Item($row, 5).Value = $data[$i][4]
}
